$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CHUT rows (B column = "CHUT"): rows 2, 8, 14, 20
# Lookback (column G) set to 2.5 months = 75.625 days, entered/displayed as 75625
# with a thousands-separator number format (#,##0 / numFmtId 3)
$chutRows = @(2, 8, 14, 20)
foreach ($r in $chutRows) {
    $cell = $ws.Range("G$r")
    $cell.Value = 75625
    $cell.NumberFormat = "#,##0"
}

# after_delivery (column H) for the second CHUT row (row 8) goes from 7 to 0
$ws.Range("H8").Value = 0

# Update the active selection on the sheet
$ws.Range("J11").Select()
